$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.934.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -5.05%  '
$ws.Range("D3").Value = "'3.280.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -5.67%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = "'557.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.74%  '
$ws.Range("D6").Value = "'185.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.67%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = "'0.592"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.16%  '
$ws.Range("D9").Value = "'3.274.75"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.51%  '
$ws.Range("D10").Value = "'0.185"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -9.13%  '
$ws.Range("D11").Value = "'0.586"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.08%  '
$ws.Range("D12").Value = "'47.43"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -7.64%  '
$ws.Range("E13").Value = '  -6.98%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = "'8.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.60%  '
$ws.Range("B15").Value = 'BitcoinCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D15").Value = "'634.84"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.40%  '
$ws.Range("D16").Value = "'3.812.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.52%  '
$ws.Range("D17").Value = "'65.965.29"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.89%  '
$ws.Range("D18").Value = "'17.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.33%  '
$ws.Range("E19").Value = '  -3.35%  '
$ws.Range("D20").Value = "'3.281.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.42%  '
$ws.Range("D21").Value = "'11.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -8.06%  '
$ws.Range("D22").Value = "'0.905"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.29%  '
$ws.Range("D23").Value = "'18.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.97%  '
$ws.Range("D24").Value = "'107.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +8.79%  '
$ws.Range("D25").Value = "'4.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -7.02%  '
$ws.Range("D26").Value = "'3.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -7.68%  '
$ws.Range("E27").Value = '  -7.30%  '
$ws.Range("D28").Value = "'9.56"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.53%  '
$ws.Range("D29").Value = "'8.68"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -6.82%  '
$ws.Range("D30").Value = "'30.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -6.89%  '
$ws.Range("D31").Value = "'3.99"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.22%  '
$ws.Range("D32").Value = "'6.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.07%  '
$ws.Range("D33").Value = "'11.03"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.91%  '
$ws.Range("E34").Value = '  -3.92%  '
$ws.Range("D35").Value = "'57.62"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.27%  '
$ws.Range("B36").Value = 'Dai'
$ws.Range("C36").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D36").Value = "'1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.04%  '
$ws.Range("B37").Value = 'Bittensor'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D37").Value = "'525.01"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.37%  '
$ws.Range("D38").Value = "'3.699.16"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.39%  '
$ws.Range("D39").Value = "'3.36"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.02%  '
$ws.Range("D40").Value = "'0.0₃0726"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -7.98%  '
$ws.Range("D41").Value = "'0.129"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.69%  '
$ws.Range("E42").Value = '  -7.56%  '
$ws.Range("D43").Value = "'32.85"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.17%  '
$ws.Range("B44").Value = 'TheGraph'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D44").Value = "'0.338"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -9.70%  '
$ws.Range("B45").Value = 'CoreDAO'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range("D45").Value = "'3.24"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -7.99%  '
$ws.Range("E46").Value = '  -1.63%  '
$ws.Range("D47").Value = "'0.0413"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.48%  '
$ws.Range("E48").Value = '  -4.12%  '
$ws.Range("E49").Value = '  -8.21%  '
$ws.Range("E50").Value = '  +0.05%  '
$ws.Range("E51").Value = '  +3.80%  '
